# Edit script: applies the interview-notes additions described in the diff.
# Hunk 1: insert 5 new sub-bullets (SignalR/WebSocket notes) after the
#         "SignalR" paragraph and before "Inlogsysteem".
# Hunk 2: trim the trailing space on "Gebruiker", insert a new bullet
#         ("Kijken of je al een paar rollen kan standaard kan maken.")
#         carrying the _GoBack bookmark, ahead of the blank paragraph and
#         the "Geel heeft een hoge prioriteit" paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Hunk 1 - insert new paragraphs after the "SignalR" bullet
# ---------------------------------------------------------------------
$signalRPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd("`r","`a") -eq "SignalR") {
        $signalRPara = $d.Paragraphs($i)
        break
    }
}
$signalRPara.Range.InsertParagraphAfter()

$d = $word.ActiveDocument
$insertedIndex = $signalRPara.Index + 1
$newPara = $d.Paragraphs($insertedIndex)
$hunk1Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">De browser hoeft niet steeds opnieuw </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>request</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> te sturen, maar wacht gewoon tot iets binnen komt.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Heeft meerdere </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fallback</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> methodes als iets niet wordt ondersteund.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Websockets</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Hub </w:t></w:r><w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>SignalR</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:t xml:space="preserve"> Client</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Nadeel van </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>front end</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> dat je </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>jqueries</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> nodig hebt.</w:t></w:r></w:p>'
$newPara.Range.InsertXML($hunk1Xml)

# ---------------------------------------------------------------------
# Hunk 2 - rework the "Gebruiker" / bookmark / "Geel heeft" paragraphs
# ---------------------------------------------------------------------
$d = $word.ActiveDocument

# Remove the existing _GoBack bookmark (it gets reinserted further down).
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$d = $word.ActiveDocument
$gebruikerPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd("`r","`a") -eq "Gebruiker ") {
        $gebruikerPara = $d.Paragraphs($i)
        break
    }
}

$geelPara = $null
for ($i = $gebruikerPara.Index; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("Geel heeft een hoge prioriteit")) {
        $geelPara = $d.Paragraphs($i)
        break
    }
}

$wholeRange = $d.Range($gebruikerPara.Range.Start, $geelPara.Range.End)
$hunk2Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Gebruiker</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Kijken of je al een paar rollen kan standaard kan maken.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p/><w:p><w:r w:rsidRPr="00EC1AE2"><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Geel heeft een hoge prioriteit</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>'
$wholeRange.InsertXML($hunk2Xml)

Write-Output "Edit complete"
